# Updated cryptos list on Sat Nov 18 15:11:19 UTC 2023 with GitHub Actions
#
# Price values in column D are stored as plain text in the source sheet
# (even when they look like numbers, e.g. "243.33"). A leading apostrophe
# forces Excel/COM to keep the assigned value as text instead of silently
# auto-converting it to a numeric cell, matching the workbook's original
# data typing.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'36.599.68"
$ws.Range("E2").Value = "  +1.37%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.942.59"
$ws.Range("E3").Value = "  +0.29%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'243.62"
$ws.Range("E5").Value = "  +0.81%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.613"
$ws.Range("E6").Value = "  +1.47%  "

# Row 7 - was USDC, now Solana
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.85"
$ws.Range("E7").Value = "  +3.20%  "

# Row 8 - was Solana, now USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.03%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.365"
$ws.Range("E9").Value = "  +0.43%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0805"
$ws.Range("E10").Value = "  -2.15%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -0.05%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "'21.97"
$ws.Range("E12").Value = "  +4.85%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "'2.226.06"
$ws.Range("E13").Value = "  +0.46%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.809"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15 - Chainlink
$ws.Range("D15").Value = "'13.34"
$ws.Range("E15").Value = "  +0.67%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "'5.18"
$ws.Range("E16").Value = "  -0.67%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "'1.944.94"
$ws.Range("E17").Value = "  -0.76%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'36.489.65"
$ws.Range("E18").Value = "  +1.36%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "'69.33"
$ws.Range("E19").Value = "  -0.25%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "'0.0₃0855"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'227.81"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.98"
$ws.Range("E22").Value = "  +0.14%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.22%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -1.85%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  +1.75%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = "  -1.34%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'159.71"
$ws.Range("E27").Value = "  -2.17%  "

# Row 28 - Kaspa
$ws.Range("D28").Value = "'0.134"
$ws.Range("E28").Value = "  +13.30%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "'19.24"
$ws.Range("E29").Value = "  +0.02%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  +1.37%  "

# Row 31 - was ImmutableX, now Filecoin
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.65"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32 - was Filecoin, now ImmutableX
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'1.09"
$ws.Range("E32").Value = "  -3.43%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.0617"

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "'4.17"
$ws.Range("E34").Value = "  -2.22%  "

# Row 35 - THORChain
$ws.Range("D35").Value = "'6.16"
$ws.Range("E35").Value = "  +2.10%  "

# Row 36 - BinanceUSD
$ws.Range("E36").Value = "  -0.11%  "

# Row 37 - was WEMIXToken, now LidoDAOToken
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.21"
$ws.Range("E37").Value = "  +2.59%  "

# Row 38 - was LidoDAOToken, now WEMIXToken
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = "  -0.82%  "

# Row 39 - RenderToken
$ws.Range("D39").Value = "'3.30"
$ws.Range("E39").Value = "  +15.38%  "

# Row 40 - Cronos
$ws.Range("D40").Value = "'0.0987"
$ws.Range("E40").Value = "  +2.46%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  +1.30%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  +0.87%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  -2.18%  "

# Row 44 - InjectiveProtocol
$ws.Range("D44").Value = "'15.81"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45 - Maker
$ws.Range("D45").Value = "'1.341.35"
$ws.Range("E45").Value = "  +0.38%  "

# Row 46 - ARBITRUM
$ws.Range("D46").Value = "'1.03"
$ws.Range("E46").Value = "  +0.05%  "

# Row 47 - Aave
$ws.Range("D47").Value = "'86.50"
$ws.Range("E47").Value = "  -1.01%  "

# Row 48 - FraxShare
$ws.Range("D48").Value = "'7.13"
$ws.Range("E48").Value = "  -2.19%  "

# Row 49 - MXToken
$ws.Range("E49").Value = "  +0.54%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "'2.119.12"
$ws.Range("E50").Value = "  +0.53%  "

# Row 51 - MultiversX
$ws.Range("D51").Value = "'43.07"
$ws.Range("E51").Value = "  -5.93%  "
